# Refresh the crypto ranking table (Sheet1) with the latest scraped values.
# Mirrors the periodic "Updated cryptos list ... with GitHub Actions" commit:
# most rows just get new Price (D) / Volume(1h) (E) figures, a few rows also
# change rank and therefore swap their Coin name + Link (B/C) between rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "52.253.93"
$ws.Range("E2").Value = "  +1.31%  "

# Row 3 (Ethereum)
$ws.Range("D3").Value = "2.829.53"
$ws.Range("E3").Value = "  +3.37%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 (BNB)
$ws.Range("D5").Value = "'356.15"
$ws.Range("E5").Value = "  +7.07%  "

# Row 6 (Solana)
$ws.Range("D6").Value = "'114.23"
$ws.Range("E6").Value = "  -1.73%  "

# Row 7 (XRP)
$ws.Range("D7").Value = "'0.548"
$ws.Range("E7").Value = "  +2.90%  "

# Row 8 (USDC)
$ws.Range("E8").Value = "  +0.03%  "

# Row 9 (Cardano)
$ws.Range("D9").Value = "'0.605"
$ws.Range("E9").Value = "  +6.06%  "

# Row 10 (Avalanche)
$ws.Range("D10").Value = "'41.88"
$ws.Range("E10").Value = "  +1.34%  "

# Row 11 (Dogecoin)
$ws.Range("D11").Value = "'0.0850"
$ws.Range("E11").Value = "  +0.06%  "

# Row 12 (Chainlink)
$ws.Range("D12").Value = "'20.10"
$ws.Range("E12").Value = "  +0.37%  "

# Row 13 (TRON)
$ws.Range("E13").Value = "  +1.42%  "

# Row 14 (Polkadot)
$ws.Range("D14").Value = "'7.76"
$ws.Range("E14").Value = "  +2.84%  "

# Row 15 (WrappedliquidstakedEther2.0)
$ws.Range("D15").Value = "3.273.38"
$ws.Range("E15").Value = "  +3.25%  "

# Row 16 (WrappedEther)
$ws.Range("D16").Value = "2.831.51"
$ws.Range("E16").Value = "  +3.10%  "

# Row 17 (Polygon)
$ws.Range("D17").Value = "'0.890"
$ws.Range("E17").Value = "  +1.60%  "

# Row 18 (WrappedBTC)
$ws.Range("D18").Value = "52.240.24"
$ws.Range("E18").Value = "  +1.40%  "

# Row 19 (ImmutableX)
$ws.Range("E19").Value = "  +2.13%  "

# Row 20 (Uniswap->InternetComputer(DFINITY))
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").Value = "'13.80"
$ws.Range("E20").Value = "  +2.81%  "

# Row 21 (InternetComputer(DFINITY)->Uniswap)
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'7.28"
$ws.Range("E21").Value = "  +7.13%  "

# Row 22 (ShibaInu)
$ws.Range("D22").Value = "0.0₃0994"
$ws.Range("E22").Value = "  +2.87%  "

# Row 23 (BitcoinCash)
$ws.Range("D23").Value = "'270.69"
$ws.Range("E23").Value = "  -3.23%  "

# Row 24 (Litecoin)
$ws.Range("D24").Value = "'69.68"
$ws.Range("E24").Value = "  +0.22%  "

# Row 25 (PancakeSwap)
$ws.Range("D25").Value = "'2.80"
$ws.Range("E25").Value = "  +6.20%  "

# Row 26 (EthereumClassic)
$ws.Range("D26").Value = "'26.74"
$ws.Range("E26").Value = "  +0.50%  "

# Row 27 (Dai)
$ws.Range("E27").Value = "  +0.10%  "

# Row 28 (Cosmos)
$ws.Range("E28").Value = "  +1.23%  "

# Row 29 (Toncoin)
$ws.Range("E29").Value = "  +1.49%  "

# Row 30 (Kaspa)
$ws.Range("D30").Value = "'0.141"
$ws.Range("E30").Value = "  +1.31%  "

# Row 31 (OKB)
$ws.Range("D31").Value = "'50.53"
$ws.Range("E31").Value = "  +0.79%  "

# Row 32 (InjectiveProtocol)
$ws.Range("D32").Value = "'33.91"
$ws.Range("E32").Value = "  -2.91%  "

# Row 33 (Filecoin)
$ws.Range("D33").Value = "'5.88"
$ws.Range("E33").Value = "  +6.26%  "

# Row 34 (VeChain)
$ws.Range("E34").Value = "  +27.71%  "

# Row 35 (Hedera)
$ws.Range("E35").Value = "  +1.94%  "

# Row 36 (FirstDigitalUSD)
$ws.Range("E36").Value = "  -0.07%  "

# Row 37 (ARBITRUM)
$ws.Range("D37").Value = "'2.09"

# Row 38 (Celestia->RenderToken)
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'4.91"
$ws.Range("E38").Value = "  -1.76%  "

# Row 39 (RenderToken->Celestia)
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").Value = "'18.52"
$ws.Range("E39").Value = "  -2.10%  "

# Row 40 (LidoDAOToken)
$ws.Range("E40").Value = "  +2.46%  "

# Row 41 (EnergySwap->Stacks)
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'2.58"
$ws.Range("E41").Value = "  +8.77%  "

# Row 42 (Stacks->Monero)
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "'128.09"
$ws.Range("E42").Value = "  +0.11%  "

# Row 43 (Monero->EnergySwap)
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'23.36"
$ws.Range("E43").Value = "  +1.78%  "

# Row 44 (Stellar)
$ws.Range("E44").Value = "  +1.94%  "

# Row 45 (WEMIXToken)
$ws.Range("E45").Value = "  +2.55%  "

# Row 46 (NEARProtocol)
$ws.Range("E46").Value = "  +1.99%  "

# Row 47 (Maker)
$ws.Range("D47").Value = "2.041.66"
$ws.Range("E47").Value = "  -2.19%  "

# Row 48 (ApeXProtocol)
$ws.Range("E48").Value = "  +3.07%  "

# Row 49 (SEI)
$ws.Range("D49").Value = "'0.975"
$ws.Range("E49").Value = "  +13.94%  "

# Row 50 (THORChain)
$ws.Range("D50").Value = "'5.71"
$ws.Range("E50").Value = "  +3.87%  "

# Row 51 (MultiversX)
$ws.Range("D51").Value = "'60.35"
$ws.Range("E51").Value = "  +0.93%  "
